{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) inside specific bullet /\n// impact paragraphs, matching the author's commit.\n//\n// Strategy: each target paragraph is located by its exact (pre-edit) text\n// so we never accidentally touch a different paragraph that merely shares\n// a substring (several numbers like \"23%\", \"87%\", \"73.5%\", \"$4.7M\" recur\n// throughout the resume). Once the paragraph is located, the metric\n// sub-strings inside it are located with a paragraph-scoped `search()`\n// (left-to-right, matching document order) and given bold + the metric\n// color (#2C3E50) \u2014 Word automatically splits the run and preserves\n// surrounding whitespace (`xml:space=\"preserve\"`) for us.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Find the single paragraph under `body` whose full text equals `text`\n// exactly (throws if zero or more than one match \u2014 keeps us honest).\nasync function getParagraphByExactText(body, text) {\n  const paras = body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n\n  const matches = paras.items.filter((p) => p.text === text);\n  if (matches.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 paragraph matching \" + JSON.stringify(text) +\n      \" but found \" + matches.length\n    );\n  }\n  return matches[0];\n}\n\n// Bold + color every metric substring found in `metrics` (in order) inside\n// `paragraph`. Each metric is searched for independently immediately\n// before it is styled, so earlier edits (which split runs) can't shift\n// the offsets of later matches.\nasync function highlightMetrics(paragraph, metrics) {\n  for (const metric of metrics) {\n    const results = paragraph.search(metric, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length < 1) {\n      throw new Error(\"Metric \" + JSON.stringify(metric) + \" not found in paragraph\");\n    }\n\n    const target = results.items[0];\n    target.font.bold = true;\n    target.font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n\nconst body = context.document.body;\n\n// 1) Professional Experience \u2014 Siege Analytics \u2014 race coding bullet.\nconst p1 = await getParagraphByExactText(\n  body,\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n);\nawait highlightMetrics(p1, [\"23%\", \"64%\"]);\n\n// 2) Professional Experience \u2014 Siege Analytics \u2014 turnout prediction bullet.\nconst p2 = await getParagraphByExactText(\n  body,\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n);\nawait highlightMetrics(p2, [\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\"]);\n\n// 3) Professional Experience \u2014 Myers Research \u2014 RFP bullet.\nconst p3 = await getParagraphByExactText(\n  body,\n  \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n);\nawait highlightMetrics(p3, [\"1,200\"]);\n\n// 4) Professional Experience \u2014 Lake Research Partners \u2014 meta-analysis bullet.\nconst p4 = await getParagraphByExactText(\n  body,\n  \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\"\n);\nawait highlightMetrics(p4, [\"$400M\", \"$1B\"]);\n\n// 5) Key Achievements and Impact \u2014 mapping cost algorithm bullet.\nconst p5 = await getParagraphByExactText(\n  body,\n  \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\"\n);\nawait highlightMetrics(p5, [\"73.5%\", \"$4.7M\"]);\n\n// 6) Key Achievements and Impact \u2014 turnout prediction bullet (short form).\nconst p6 = await getParagraphByExactText(\n  body,\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n);\nawait highlightMetrics(p6, [\"87%\", \"71%\"]);\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) inside specific bullet /\n# impact paragraphs, matching the author's commit.\n#\n# Strategy: each target paragraph is located by its exact (pre-edit) text\n# (iterating $d.Paragraphs and comparing Range.Text with the trailing\n# paragraph mark trimmed off) so we never accidentally touch a different\n# paragraph that merely shares a substring (several numbers like \"23%\",\n# \"87%\", \"73.5%\", \"$4.7M\" recur throughout the resume). Once the paragraph\n# is located, the metric sub-strings inside it are located with a\n# paragraph-scoped Find.Execute (left-to-right, matching document order)\n# and given bold + the metric color (#2C3E50) -- Word automatically splits\n# the run and preserves surrounding whitespace (xml:space=\"preserve\") for\n# us.\n\n$HighlightColorHex = \"2C3E50\"\n\n# RGB hex -> Word's OLE/wdColor integer (0x00BBGGRR).\nfunction Get-WdColor {\n    param([string]$hex)\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return $r + ($g * 256) + ($b * 65536)\n}\n\n# Find the single paragraph under $doc whose full text equals $text exactly\n# (throws if zero or more than one match -- keeps us honest).\nfunction Get-ParagraphByExactText {\n    param($doc, $text)\n    $found = @()\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]0x0D, [char]0x07)\n        if ($t -eq $text) {\n            $found += $p\n        }\n    }\n    if ($found.Count -ne 1) {\n        throw (\"Expected exactly 1 paragraph matching \" + $text + \" but found \" + $found.Count)\n    }\n    return $found[0]\n}\n\n# Bold + color the first occurrence of $metricText inside $para (scoped\n# Find.Execute, so an identical number elsewhere in the document is never\n# touched).\nfunction Set-MetricHighlight {\n    param($para, $metricText)\n    $rng = $para.Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $metricText\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $wasFound = $rng.Find.Execute()\n    if (-not $wasFound) {\n        throw (\"Metric \" + $metricText + \" not found in paragraph\")\n    }\n    $rng.Font.Bold = 1\n    $wdColor = Get-WdColor $HighlightColorHex\n    $rng.Font.Color = $wdColor\n}\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n$plusMinus = [char]0xB1\n\n# 1) Professional Experience -- Siege Analytics -- race coding bullet.\n$text1 = $bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n$p1 = Get-ParagraphByExactText $d $text1\nSet-MetricHighlight $p1 \"23%\"\nSet-MetricHighlight $p1 \"64%\"\n\n# 2) Professional Experience -- Siege Analytics -- turnout prediction bullet.\n$m2c = $plusMinus + \"4.2%\"\n$m2d = $plusMinus + \"2.1%\"\n$text2 = $bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + $m2c + \" to \" + $m2d\n$p2 = Get-ParagraphByExactText $d $text2\nSet-MetricHighlight $p2 \"87%\"\nSet-MetricHighlight $p2 \"71%\"\nSet-MetricHighlight $p2 $m2c\nSet-MetricHighlight $p2 $m2d\n\n# 3) Professional Experience -- Myers Research -- RFP bullet.\n$text3 = $bullet + \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n$p3 = Get-ParagraphByExactText $d $text3\nSet-MetricHighlight $p3 \"1,200\"\n\n# 4) Professional Experience -- Lake Research Partners -- meta-analysis bullet.\n$text4 = $bullet + \" Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n$p4 = Get-ParagraphByExactText $d $text4\nSet-MetricHighlight $p4 \"`$400M\"\nSet-MetricHighlight $p4 \"`$1B\"\n\n# 5) Key Achievements and Impact -- mapping cost algorithm bullet.\n$text5 = $bullet + \" Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n$p5 = Get-ParagraphByExactText $d $text5\nSet-MetricHighlight $p5 \"73.5%\"\nSet-MetricHighlight $p5 \"`$4.7M\"\n\n# 6) Key Achievements and Impact -- turnout prediction bullet (short form).\n$text6 = $bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n$p6 = Get-ParagraphByExactText $d $text6\nSet-MetricHighlight $p6 \"87%\"\nSet-MetricHighlight $p6 \"71%\"\n"}
